$wb = $excel.ActiveWorkbook

# --- Content changes: rename "Logical"/"logical" -> "Conceptual"/"conceptual" ---

# Properties sheet: header row column O ("Logical" -> "Conceptual")
$wsProps = $wb.Worksheets.Item("Properties")
$wsProps.Range("O2").Value = "Conceptual"

# Views sheet: header row column G ("Logical" -> "Conceptual")
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Range("G2").Value = "Conceptual"

# Metadata sheet: row 11 value ("logical" -> "conceptual")
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A11").Value = "conceptual"

# --- View/selection changes ---

# Narrow the Properties sheet selection from O3:O7 down to O3
$wsProps.Range("O3").Select()

# Narrow the Views sheet selection from G3:G5 down to G3
$wsViews.Range("G3").Select()

# Make Metadata the active/tab-selected sheet with selection at A12
$wsMeta.Activate()
$wsMeta.Range("A12").Select()
